$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 6
$ws.Range("C4").Value = 6
$ws.Range("C5").Value = 8
$ws.Range("C6").Value = 7
$ws.Range("C7").Value = 6
$ws.Range("C8").Value = 6
$ws.Range("B9").Value = "<you>"
$ws.Range("C9").Value = 8
$ws.Range("B10").Value = "<then>"
$ws.Range("C10").Value = 7
$ws.Range("C11").Value = 11
$ws.Range("C12").Value = 6
$ws.Range("C13").Value = 9
$ws.Range("C14").Value = 4
$ws.Range("C15").Value = 5
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 7
$ws.Range("B18").Value = "<of>"
$ws.Range("C18").Value = 9
$ws.Range("C19").Value = 6
$ws.Range("C20").Value = 4
$ws.Range("C21").Value = 7
$ws.Range("C22").Value = 3
$ws.Range("C23").Value = 4
$ws.Range("B24").Value = "<their>"
$ws.Range("C24").Value = 6
$ws.Range("C26").Value = 6
$ws.Range("C28").Value = 5
$ws.Range("C30").Value = 13
$ws.Range("C32").Value = 7
$ws.Range("C34").Value = 8
$ws.Range("C35").Value = 3
$ws.Range("B36").Value = "<its>"
$ws.Range("C36").Value = 8
$ws.Range("C37").Value = 7
$ws.Range("C38").Value = 6
$ws.Range("C39").Value = 2
$ws.Range("C40").Value = 6
$ws.Range("C41").Value = 5
$ws.Range("C42").Value = 6
$ws.Range("C43").Value = 8
$ws.Range("C44").Value = 4
$ws.Range("C45").Value = 11
$ws.Range("B46").Value = "<their>"
$ws.Range("C46").Value = 7
$ws.Range("C47").Value = 7
$ws.Range("B48").Value = "<up>"
$ws.Range("C48").Value = 8
$ws.Range("C49").Value = 3
$ws.Range("B50").Value = "<xirst>"
$ws.Range("C50").Value = 10
$ws.Range("C52").Value = 2
